# Apply the changes described by the commit:
#  "adjusted readme, added eventMember and item"
#
# 1. Add a new "item" worksheet (app_items table DDL) after "address".
# 2. Add a new "events_members" worksheet (app_events_members table DDL)
#    after "item", and make it the active/selected sheet.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# New sheet: item
# ---------------------------------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$itemSheet = $wb.Worksheets.Add($null, $lastSheet)
$itemSheet.Name = "item"

$itemSheet.Range("A1").Value = "create table app_items("
$itemSheet.Range("B2").Value = "item_id serial,"
$itemSheet.Range("B3").Value = "item varchar(20) not null,"
$itemSheet.Range("B4").Value = "comment varchar(1000) not null,"
$itemSheet.Range("B5").Value = "event_id integer not null,"
$itemSheet.Range("B6").Value = "member_id integer not null,"
$itemSheet.Range("B7").Value = "constraint app_item_id_pk primary key (item_id),"
$itemSheet.Range("B8").Value = "constraint app_event_fk foreign key (event_id) references app_events,"
$itemSheet.Range("B9").Value = "constraint app_member_id_fk foreign key (member_id) references app_members"
$itemSheet.Range("A10").Value = ");"

[void]$itemSheet.Range("C12").Select()

# ---------------------------------------------------------------------
# New sheet: events_members
# ---------------------------------------------------------------------
$lastSheet2 = $wb.Worksheets.Item($wb.Worksheets.Count)
$eventsMembersSheet = $wb.Worksheets.Add($null, $lastSheet2)
$eventsMembersSheet.Name = "events_members"

$eventsMembersSheet.Range("A1").Value = "create table app_events_members("
$eventsMembersSheet.Range("B2").Value = "event_id integer,"
$eventsMembersSheet.Range("B3").Value = "member_id integer,"
$eventsMembersSheet.Range("B4").Value = "constraint app_event_member_pk primary key (event_id, member_id),"
$eventsMembersSheet.Range("B5").Value = "constraint app_event_id_fk foreign key (event_id) references app_events,"
$eventsMembersSheet.Range("B6").Value = "constraint app_member_id_fk foreign key (member_id) references app_members"
$eventsMembersSheet.Range("A7").Value = ");"

[void]$eventsMembersSheet.Range("B2:B3").Select()
